$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.122044
$ws.Range("H2").Value = 15.366132
$ws.Range("I2").Value = 0.0387196063811631
$ws.Range("J2").Value = 0.0387196063811631
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.438907666666667
$ws.Range("N2").Value = 10.316723
$ws.Range("O2").Value = 0.05825422340060618
$ws.Range("P2").Value = 0.05825422340060618
$ws.Range("Q2").Value = 17.614236380604
$ws.Range("R2").Value = 158.528127425436
$ws.Range("S2").Value = 0.002255580600111812
$ws.Range("T2").Value = 0.002255580600111812
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.122044
$ws.Range("H3").Value = 15.366132
$ws.Range("I3").Value = 0.0387196063811631
$ws.Range("J3").Value = 0.0387196063811631
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.383857
$ws.Range("N3").Value = 31.151571
$ws.Range("O3").Value = 0.1758999031294962
$ws.Range("P3").Value = 0.1758999031294962
$ws.Range("Q3").Value = 53.186572443708
$ws.Range("R3").Value = 478.679151993372
$ws.Range("S3").Value = 0.00681077501165881
$ws.Range("T3").Value = 0.00681077501165881
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.122044
$ws.Range("H4").Value = 15.366132
$ws.Range("I4").Value = 0.0387196063811631
$ws.Range("J4").Value = 0.0387196063811631
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.34077833333333
$ws.Range("N4").Value = 76.022335
$ws.Range("O4").Value = 0.4292663558501786
$ws.Range("P4").Value = 0.4292663558501786
$ws.Range("Q4").Value = 129.79658161758
$ws.Range("R4").Value = 1168.16923455822
$ws.Range("S4").Value = 0.0166210243311952
$ws.Range("T4").Value = 0.0166210243311952
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.122044
$ws.Range("H5").Value = 15.366132
$ws.Range("I5").Value = 0.0387196063811631
$ws.Range("J5").Value = 0.0387196063811631
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.86921833333334
$ws.Range("N5").Value = 59.60765500000001
$ws.Range("O5").Value = 0.336579517619719
$ws.Range("P5").Value = 0.336579517619719
$ws.Range("Q5").Value = 101.77101054894
$ws.Range("R5").Value = 915.9390949404601
$ws.Range("S5").Value = 0.01303222643819727
$ws.Range("T5").Value = 0.01303222643819727
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 75.68093133333333
$ws.Range("H6").Value = 227.042794
$ws.Range("I6").Value = 0.5721028307813247
$ws.Range("J6").Value = 0.5721028307813247
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.438907666666667
$ws.Range("N6").Value = 10.316723
$ws.Range("O6").Value = 0.05825422340060618
$ws.Range("P6").Value = 0.05825422340060618
$ws.Range("Q6").Value = 260.2597349826735
$ws.Range("R6").Value = 2342.337614844062
$ws.Range("S6").Value = 0.03332740611245449
$ws.Range("T6").Value = 0.03332740611245449
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 75.68093133333333
$ws.Range("H7").Value = 227.042794
$ws.Range("I7").Value = 0.5721028307813247
$ws.Range("J7").Value = 0.5721028307813247
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.383857
$ws.Range("N7").Value = 31.151571
$ws.Range("O7").Value = 0.1758999031294962
$ws.Range("P7").Value = 0.1758999031294962
$ws.Range("Q7").Value = 785.8599685921528
$ws.Range("R7").Value = 7072.739717329375
$ws.Range("S7").Value = 0.1006328325145456
$ws.Range("T7").Value = 0.1006328325145456
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 75.68093133333333
$ws.Range("H8").Value = 227.042794
$ws.Range("I8").Value = 0.5721028307813247
$ws.Range("J8").Value = 0.5721028307813247
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.34077833333333
$ws.Range("N8").Value = 76.022335
$ws.Range("O8").Value = 0.4292663558501786
$ws.Range("P8").Value = 0.4292663558501786
$ws.Range("Q8").Value = 1917.813704978221
$ws.Range("R8").Value = 17260.32334480399
$ws.Range("S8").Value = 0.2455844973410707
$ws.Range("T8").Value = 0.2455844973410707
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 75.68093133333333
$ws.Range("H9").Value = 227.042794
$ws.Range("I9").Value = 0.5721028307813247
$ws.Range("J9").Value = 0.5721028307813247
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.86921833333334
$ws.Range("N9").Value = 59.60765500000001
$ws.Range("O9").Value = 0.336579517619719
$ws.Range("P9").Value = 0.336579517619719
$ws.Range("Q9").Value = 1503.720948332008
$ws.Range("R9").Value = 13533.48853498807
$ws.Range("S9").Value = 0.192558094813254
$ws.Range("T9").Value = 0.192558094813254
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 41.00894566666667
$ws.Range("H10").Value = 123.026837
$ws.Range("I10").Value = 0.3100032397847104
$ws.Range("J10").Value = 0.3100032397847104
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.438907666666667
$ws.Range("N10").Value = 10.316723
$ws.Range("O10").Value = 0.05825422340060618
$ws.Range("P10").Value = 0.05825422340060618
$ws.Range("Q10").Value = 141.0259776550168
$ws.Range("R10").Value = 1269.233798895151
$ws.Range("S10").Value = 0.01805899798533021
$ws.Range("T10").Value = 0.01805899798533021
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 41.00894566666667
$ws.Range("H11").Value = 123.026837
$ws.Range("I11").Value = 0.3100032397847104
$ws.Range("J11").Value = 0.3100032397847104
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 10.383857
$ws.Range("N11").Value = 31.151571
$ws.Range("O11").Value = 0.1758999031294962
$ws.Range("P11").Value = 0.1758999031294962
$ws.Range("Q11").Value = 425.8310275234364
$ws.Range("R11").Value = 3832.479247710927
$ws.Range("S11").Value = 0.05452953984796054
$ws.Range("T11").Value = 0.05452953984796054
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 41.00894566666667
$ws.Range("H12").Value = 123.026837
$ws.Range("I12").Value = 0.3100032397847104
$ws.Range("J12").Value = 0.3100032397847104
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 25.34077833333333
$ws.Range("N12").Value = 76.022335
$ws.Range("O12").Value = 0.4292663558501786
$ws.Range("P12").Value = 0.4292663558501786
$ws.Range("Q12").Value = 1039.198601822711
$ws.Range("R12").Value = 9352.787416404395
$ws.Range("S12").Value = 0.1330739610441317
$ws.Range("T12").Value = 0.1330739610441318
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 41.00894566666667
$ws.Range("H13").Value = 123.026837
$ws.Range("I13").Value = 0.3100032397847104
$ws.Range("J13").Value = 0.3100032397847104
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.86921833333334
$ws.Range("N13").Value = 59.60765500000001
$ws.Range("O13").Value = 0.336579517619719
$ws.Range("P13").Value = 0.336579517619719
$ws.Range("Q13").Value = 814.8156950708041
$ws.Range("R13").Value = 7333.341255637236
$ws.Range("S13").Value = 0.1043407409072879
$ws.Range("T13").Value = 0.1043407409072879
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 10.473618
$ws.Range("H14").Value = 31.420854
$ws.Range("I14").Value = 0.07917432305280171
$ws.Range("J14").Value = 0.0791743230528017
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.438907666666667
$ws.Range("N14").Value = 10.316723
$ws.Range("O14").Value = 0.05825422340060618
$ws.Range("P14").Value = 0.05825422340060618
$ws.Range("Q14").Value = 36.017805237938
$ws.Range("R14").Value = 324.160247141442
$ws.Range("S14").Value = 0.004612238702709675
$ws.Range("T14").Value = 0.004612238702709674
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 10.473618
$ws.Range("H15").Value = 31.420854
$ws.Range("I15").Value = 0.07917432305280171
$ws.Range("J15").Value = 0.0791743230528017
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 10.383857
$ws.Range("N15").Value = 31.151571
$ws.Range("O15").Value = 0.1758999031294962
$ws.Range("P15").Value = 0.1758999031294962
$ws.Range("Q15").Value = 108.756551584626
$ws.Range("R15").Value = 978.8089642616339
$ws.Range("S15").Value = 0.01392675575533126
$ws.Range("T15").Value = 0.01392675575533125
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 10.473618
$ws.Range("H16").Value = 31.420854
$ws.Range("I16").Value = 0.07917432305280171
$ws.Range("J16").Value = 0.0791743230528017
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 25.34077833333333
$ws.Range("N16").Value = 76.022335
$ws.Range("O16").Value = 0.4292663558501786
$ws.Range("P16").Value = 0.4292663558501786
$ws.Range("Q16").Value = 265.40963208601
$ws.Range("R16").Value = 2388.68668877409
$ws.Range("S16").Value = 0.03398687313378098
$ws.Range("T16").Value = 0.03398687313378097
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 10.473618
$ws.Range("H17").Value = 31.420854
$ws.Range("I17").Value = 0.07917432305280171
$ws.Range("J17").Value = 0.0791743230528017
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.86921833333334
$ws.Range("N17").Value = 59.60765500000001
$ws.Range("O17").Value = 0.336579517619719
$ws.Range("P17").Value = 0.336579517619719
$ws.Range("Q17").Value = 208.10260278193
$ws.Range("R17").Value = 1872.92342503737
$ws.Range("S17").Value = 0.0266484554609798
$ws.Range("T17").Value = 0.0266484554609798
